$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.834.30'
$ws.Range('E2').Value = '  -1.69%  '
$ws.Range('D3').Value = '3.512.43'
$ws.Range('E3').Value = '  -3.47%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '578.76'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -4.33%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '192.85'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -4.01%  '
$ws.Range('D7').Value = '3.503.82'
$ws.Range('E7').Value = '  -3.34%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.613'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.46%  '
$ws.Range('E9').Value = '  +0.10%  '
$ws.Range('E10').Value = '  -6.76%  '
$ws.Range('E11').Value = '  -4.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '51.53'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.28%  '
$ws.Range('E13').Value = '  -6.46%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '9.19'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.28%  '
$ws.Range('D15').Value = '4.072.10'
$ws.Range('E15').Value = '  -3.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '647.54'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -4.96%  '
$ws.Range('D17').Value = '69.820.04'
$ws.Range('E17').Value = '  -1.73%  '
$ws.Range('D18').Value = '3.517.21'
$ws.Range('E18').Value = '  -2.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.72%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.42'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.45%  '
$ws.Range('B21').Value = 'TRON'
$ws.Range('C21').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.121'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.952'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.87%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '18.16'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.36'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.62%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '99.55'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -5.47%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.30'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -7.16%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.90'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.14'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.07%  '
$ws.Range('E29').Value = '  -5.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '32.79'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.61%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -7.22%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.77'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.62%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.68'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -4.42%  '
$ws.Range('B34').Value = 'Bittensor'
$ws.Range('C34').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '573.78'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +10.54%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.110'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -4.70%  '
$ws.Range('D37').Value = '3.775.37'
$ws.Range('E37').Value = '  -3.91%  '
$ws.Range('B38').Value = 'Dai'
$ws.Range('C38').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.00'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.24%  '
$ws.Range('B39').Value = 'CoreDAO'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.86'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +45.39%  '
$ws.Range('E40').Value = '  -8.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.62'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.14%  '
$ws.Range('E42').Value = '  -4.17%  '
$ws.Range('E43').Value = '  -3.97%  '
$ws.Range('E44').Value = '  -2.37%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '34.38'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.10%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0445'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.96%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.88'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.22%  '
$ws.Range('E48').Value = '  -3.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.136'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -3.48%  '
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.22'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.08%  '
